# Update crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.331.99'
$ws.Range('E2').Value = '  +0.03%  '

$ws.Range('D3').Value = '1.711.42'
$ws.Range('E3').Value = '  +0.09%  '

$ws.Range('E4').Value = '  +0.65%  '

$ws.Range('D5').Value = "'224.37"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.21%  '

$ws.Range('D6').Value = "'0.5264"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.61%  '

$ws.Range('D7').Value = "'1.009"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.62%  '

$ws.Range('D8').Value = "'0.06624"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.11%  '

$ws.Range('D9').Value = "'0.2636"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.21%  '

$ws.Range('D10').Value = "'20.65"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.54%  '

$ws.Range('D11').Value = "'0.07758"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.59%  '

$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.725.45'
$ws.Range('E12').Value = '  +0.44%  '

$ws.Range('D13').Value = '1.948.26'
$ws.Range('E13').Value = '  +0.15%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = "'4.437"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.88%  '

$ws.Range('D15').Value = "'0.5763"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.32%  '

$ws.Range('D16').Value = '0.0₅8158'
$ws.Range('E16').Value = '  -0.46%  '

$ws.Range('D17').Value = "'67.45"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.24%  '

$ws.Range('D18').Value = '27.356.71'
$ws.Range('E18').Value = '  +0.18%  '

$ws.Range('D19').Value = "'218.16"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '

$ws.Range('D20').Value = "'1.007"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.48%  '

$ws.Range('D21').Value = "'4.630"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.10%  '

$ws.Range('D22').Value = "'10.39"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.69%  '

$ws.Range('E23').Value = '  +0.88%  '

$ws.Range('E24').Value = '  +0.61%  '

$ws.Range('D25').Value = "'145.38"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.85%  '

$ws.Range('E26').Value = '  -1.54%  '

$ws.Range('D27').Value = "'0.1198"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.75%  '

$ws.Range('D28').Value = "'7.186"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.13%  '

$ws.Range('D29').Value = "'16.12"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.32%  '

$ws.Range('D30').Value = "'0.05296"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.76%  '

$ws.Range('D31').Value = "'1.293"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.09%  '

$ws.Range('D32').Value = "'3.464"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.82%  '

$ws.Range('D33').Value = "'3.342"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.15%  '

$ws.Range('D34').Value = "'1.635"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.27%  '

$ws.Range('D35').Value = "'2.831"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.28%  '

$ws.Range('D36').Value = "'0.9477"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.44%  '

$ws.Range('D37').Value = "'2.403"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.57%  '

$ws.Range('D38').Value = "'0.5861"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.20%  '

$ws.Range('D39').Value = '1.178.57'
$ws.Range('E39').Value = '  +13.18%  '

$ws.Range('D40').Value = "'0.01647"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.11%  '

$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = "'1.009"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.61%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'5.764"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.77%  '

$ws.Range('D43').Value = "'0.8381"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.09%  '

$ws.Range('D44').Value = "'100.94"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.18%  '

$ws.Range('D45').Value = '1.856.52'
$ws.Range('E45').Value = '  +0.18%  '

$ws.Range('E46').Value = '  +2.81%  '

$ws.Range('D47').Value = "'57.36"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.25%  '

$ws.Range('D48').Value = "'0.4558"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.39%  '

$ws.Range('D49').Value = "'1.004"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.26%  '

$ws.Range('D50').Value = "'8.115"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.64%  '

$ws.Range('D51').Value = "'0.05233"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.15%  '
